$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it looks like a number
# (e.g. "1.000", "34.84", "0.000007599"), so Excel does not silently coerce
# it to a numeric type and strip the formatting that the source data relies
# on (trailing zeros, thousands "." separators, etc). We briefly mark the
# cell as Text, assign the value, then restore the default "Normal" style so
# the cell's formatting matches an untouched cell.
function Set-CellText {
    param(
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-CellText 2 4 '30.580.69'
Set-CellText 2 5 '  +1.27%  '
Set-CellText 3 4 '1.919.06'
Set-CellText 3 5 '  +2.45%  '
Set-CellText 4 4 '1.000'
Set-CellText 4 5 '  +0.69%  '
Set-CellText 5 4 '247.40'
Set-CellText 6 5 '  +0.65%  '
Set-CellText 7 4 '0.4723'
Set-CellText 7 5 '  +1.39%  '
Set-CellText 8 4 '0.2874'
Set-CellText 8 5 '  +1.80%  '
Set-CellText 9 4 '0.06839'
Set-CellText 9 5 '  +5.42%  '
Set-CellText 10 5 '  -2.36%  '
Set-CellText 11 4 '18.32'
Set-CellText 11 5 '  -1.77%  '
Set-CellText 12 4 '1.916.59'
Set-CellText 12 5 '  +2.49%  '
Set-CellText 13 4 '0.07693'
Set-CellText 13 5 '  +2.50%  '
Set-CellText 14 4 '5.282'
Set-CellText 14 5 '  +5.08%  '
Set-CellText 15 4 '0.6688'
Set-CellText 15 5 '  +5.65%  '
Set-CellText 16 4 '288.74'
Set-CellText 16 5 '  -7.17%  '
Set-CellText 17 4 '30.584.22'
Set-CellText 17 5 '  +1.39%  '
Set-CellText 18 2 'ShibaInu'
Set-CellText 18 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText 18 4 '0.000007599'
Set-CellText 18 5 '  +2.05%  '
Set-CellText 19 2 'Dai'
Set-CellText 19 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText 19 4 '1.000'
Set-CellText 19 5 '  +0.71%  '
Set-CellText 20 5 '  +1.55%  '
Set-CellText 21 4 '5.555'
Set-CellText 21 5 '  +9.99%  '
Set-CellText 22 4 '2.167.36'
Set-CellText 22 5 '  +3.79%  '
Set-CellText 23 4 '1.000'
Set-CellText 23 5 '  +0.68%  '
Set-CellText 24 4 '6.298'
Set-CellText 24 5 '  +1.74%  '
Set-CellText 25 4 '9.373'
Set-CellText 25 5 '  +2.24%  '
Set-CellText 26 4 '168.66'
Set-CellText 26 5 '  +2.41%  '
Set-CellText 27 5 '  +5.54%  '
Set-CellText 28 4 '2.120'
Set-CellText 28 5 '  +7.25%  '
Set-CellText 29 4 '0.1069'
Set-CellText 29 5 '  -0.93%  '
Set-CellText 30 5 '  +4.64%  '
Set-CellText 31 4 '4.174'
Set-CellText 31 5 '  +3.27%  '
Set-CellText 32 4 '4.086'
Set-CellText 32 5 '  +5.39%  '
Set-CellText 33 4 '0.05041'
Set-CellText 33 5 '  +2.55%  '
Set-CellText 34 5 '  -0.73%  '
Set-CellText 35 4 '1.148'
Set-CellText 35 5 '  +1.33%  '
Set-CellText 36 4 '0.02068'
Set-CellText 36 5 '  +7.35%  '
Set-CellText 37 4 '2.747'
Set-CellText 37 5 '  +1.70%  '
Set-CellText 38 4 '2.689'
Set-CellText 38 5 '  +1.10%  '
Set-CellText 39 4 '2.058'
Set-CellText 39 5 '  +3.48%  '
Set-CellText 40 4 '111.15'
Set-CellText 40 5 '  +4.00%  '
Set-CellText 41 4 '0.8802'
Set-CellText 41 5 '  +1.76%  '
Set-CellText 42 4 '0.4392'
Set-CellText 42 5 '  +7.43%  '
Set-CellText 43 4 '5.878'
Set-CellText 43 5 '  +2.52%  '
Set-CellText 44 4 '1.000'
Set-CellText 45 4 '67.20'
Set-CellText 45 5 '  +0.70%  '
Set-CellText 46 4 '7.264'
Set-CellText 46 5 '  +1.70%  '
Set-CellText 47 4 '9.269'
Set-CellText 47 5 '  +0.58%  '
Set-CellText 48 5 '  +15.32%  '
Set-CellText 49 4 '0.1232'
Set-CellText 49 5 '  +2.79%  '
Set-CellText 50 2 'Elrond'
Set-CellText 50 3 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-CellText 50 4 '34.84'
Set-CellText 50 5 '  +2.07%  '
Set-CellText 51 2 'Decentraland'
Set-CellText 51 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-CellText 51 4 '0.4067'
Set-CellText 51 5 '  +8.29%  '
